$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $value) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "324.98"
Set-TextValue "E2" "-2.00%"
Set-TextValue "D3" "44.32"
Set-TextValue "E3" "0.43%"
Set-TextValue "D4" "5.523"
Set-TextValue "E4" "-4.39%"
Set-TextValue "D5" "0.08047"
Set-TextValue "E5" "-3.54%"
Set-TextValue "D6" "8.713"
Set-TextValue "E6" "-1.15%"
Set-TextValue "D7" "4.345"
Set-TextValue "E7" "-3.48%"
Set-TextValue "D8" "1.896"
Set-TextValue "E8" "-3.33%"
Set-TextValue "D9" "2.739"
Set-TextValue "E9" "-6.26%"
Set-TextValue "E10" "1.56%"
Set-TextValue "D11" "0.1179"
Set-TextValue "E11" "-4.60%"
Set-TextValue "D12" "0.1894"
Set-TextValue "E12" "-3.13%"
Set-TextValue "D13" "0.09881"
Set-TextValue "E13" "4.17%"
Set-TextValue "D14" "0.04186"
Set-TextValue "E14" "5.80%"
Set-TextValue "D16" "0.001277"
Set-TextValue "E16" "-2.04%"
Set-TextValue "D17" "0.005998"
Set-TextValue "E17" "1.14%"
Set-TextValue "D18" "3.599"
Set-TextValue "E18" "2.68%"
Set-TextValue "D20" "8.447"
Set-TextValue "E20" "-6.51%"
Set-TextValue "D21" "0.1375"
Set-TextValue "E21" "0.24%"
Set-TextValue "D22" "0.2537"
Set-TextValue "E22" "-1.27%"
Set-TextValue "D23" "0.04256"
Set-TextValue "E23" "-3.40%"
Set-TextValue "D24" "0.001241"
Set-TextValue "E24" "-1.16%"
Set-TextValue "D25" "0.004489"
Set-TextValue "E25" "2.16%"
Set-TextValue "D26" "0.0001236"
Set-TextValue "E26" "3.82%"
Set-TextValue "D27" "0.0004006"
Set-TextValue "E27" "0.42%"
Set-TextValue "D39" "0.02641"
Set-TextValue "E39" "-6.94%"
Set-TextValue "D40" "0.05494"
Set-TextValue "E40" "-2.68%"
Set-TextValue "D41" "0.007658"
Set-TextValue "E41" "-3.08%"
Set-TextValue "D42" "0.1391"
Set-TextValue "E42" "-2.45%"
Set-TextValue "D43" "0.006746"
Set-TextValue "E43" "-25.62%"
Set-TextValue "D44" "0.002056"
Set-TextValue "E44" "-2.13%"
Set-TextValue "D45" "0.009210"
Set-TextValue "E45" "-7.29%"
Set-TextValue "D46" "0.00007149"
Set-TextValue "E46" "-1.88%"
Set-TextValue "D47" "0.00000000754"
Set-TextValue "E47" "0.45%"
Set-TextValue "D48" "0.003431"
Set-TextValue "E48" "-13.58%"
Set-TextValue "D49" "0.002280"
Set-TextValue "E49" "0.08%"
Set-TextValue "D50" "0.00002110"
Set-TextValue "E50" "0.45%"
Set-TextValue "D51" "0.0002009"
Set-TextValue "E51" "0.45%"
